# Apply edits described by the diff: update several numeric outputs,
# remove the trailing "Class I weight estimation" row (39), and add a
# new "Stability margin" row (38).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update single-value cells
$ws.Range("B8").Value = 6000

$ws.Range("B23").Value = 19235.93618374716
$ws.Range("B24").Value = 85128.60044282374
$ws.Range("B25").Value = 54568.66425907657

$ws.Range("B28").Value = 1246.114006310022
$ws.Range("B29").Value = 9694.815192239563
$ws.Range("B30").Value = 6160.906687053855
$ws.Range("B31").Value = 1337.865595885331
$ws.Range("B32").Value = 77.37596008705309
$ws.Range("B33").Value = 389.5769267962468

$ws.Range("B36").Value = 10.01505877647578
$ws.Range("B37").Value = 10.82710767237309

# Replace old row 39 ("Class I weight estimation") with the new row 38
# ("Stability margin"). Clear row 39 entirely and write new content to
# row 38.
$ws.Range("A39:C39").ClearContents()

$ws.Range("A38").Value = "Stability margin"
$ws.Range("B38").Value = -6.540066047709018
$ws.Range("C38").Value = "m"
